$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 157.1
$ws.Range("I39").Value = 80.111115
$ws.Range("K39").Value = 240.333345
$ws.Range("M39").Value = 55.66665499999999
$ws.Range("H42").Value = 762
$ws.Range("J42").Value = 24
$ws.Range("L42").Value = 72
$ws.Range("N42").Value = -532
$ws.Range("H62").Value = 7152092
$ws.Range("I62").Value = 10215474
$ws.Range("J62").Value = 4201.8335
$ws.Range("K62").Value = 10215474
$ws.Range("L62").Value = 4201.8335
$ws.Range("M62").Value = -10214850
$ws.Range("N62").Value = -5449.8335
$ws.Range("H65").Value = 7152092
$ws.Range("I65").Value = 10215474
$ws.Range("J65").Value = 4201.8335
$ws.Range("K65").Value = 51077370
$ws.Range("L65").Value = 21009.1675
$ws.Range("M65").Value = -51074250
$ws.Range("N65").Value = -27249.1675
$ws.Range("H127").Value = 7539.8335
$ws.Range("I127").Value = 6806.3125
$ws.Range("K127").Value = 20418.9375
$ws.Range("M127").Value = -15458.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2979973
$ws.Range("I134").Value = 3402629
$ws.Range("K134").Value = 10207887
$ws.Range("M134").Value = -10205352

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4500
$ws.Range("I6").Value = 4500
$ws.Range("K6").Value = 4500
$ws.Range("M6").Value = -4387
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("N56").Value = 0
$ws.Range("H58").Value = 25644062
$ws.Range("I58").Value = 35716076
$ws.Range("K58").Value = 35716076
$ws.Range("M58").Value = -35715873
$ws.Range("H136").Value = 25644062
$ws.Range("I136").Value = 35716076
$ws.Range("K136").Value = 107148228
$ws.Range("M136").Value = -107145678

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 359254.5
$ws.Range("I5").Value = 2689.8333
$ws.Range("K5").Value = 8069.499899999999
$ws.Range("M5").Value = -7957.499899999999
$ws.Range("H46").Value = 1457.8889
$ws.Range("I46").Value = 314.75
$ws.Range("K46").Value = 944.25
$ws.Range("M46").Value = -853.25
$ws.Range("H63").Value = 13471.846
$ws.Range("I63").Value = 7947.857
$ws.Range("J63").Value = 19916.5
$ws.Range("K63").Value = 23843.571
$ws.Range("L63").Value = 59749.5
$ws.Range("M63").Value = -23094.571
$ws.Range("N63").Value = -61247.5
$ws.Range("H66").Value = 13471.846
$ws.Range("I66").Value = 7947.857
$ws.Range("J66").Value = 19916.5
$ws.Range("K66").Value = 71530.713
$ws.Range("L66").Value = 179248.5
$ws.Range("M66").Value = -67786.713
$ws.Range("N66").Value = -186736.5
$ws.Range("H135").Value = 359254.5
$ws.Range("I135").Value = 2689.8333
$ws.Range("K135").Value = 24208.4997
$ws.Range("M135").Value = -21673.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4522
$ws.Range("I22").Value = 4362.6665
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 4362.6665
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -3833.6665
$ws.Range("N22").Value = -6058
$ws.Range("H41").Value = 19055.445
$ws.Range("I41").Value = 17666.334
$ws.Range("K41").Value = 17666.334
$ws.Range("M41").Value = -17311.334
$ws.Range("H70").Value = 12621
$ws.Range("J70").Value = 13530.143
$ws.Range("L70").Value = 13530.143
$ws.Range("N70").Value = -14070.143
$ws.Range("H73").Value = 12621
$ws.Range("J73").Value = 13530.143
$ws.Range("L73").Value = 13530.143
$ws.Range("N73").Value = -15402.143
$ws.Range("H80").Value = 5467.3335
$ws.Range("I80").Value = 4021.7144
$ws.Range("J80").Value = 7491.2
$ws.Range("K80").Value = 4021.7144
$ws.Range("L80").Value = 7491.2
$ws.Range("M80").Value = -3023.7144
$ws.Range("N80").Value = -9487.200000000001
$ws.Range("H83").Value = 5467.3335
$ws.Range("I83").Value = 4021.7144
$ws.Range("J83").Value = 7491.2
$ws.Range("K83").Value = 20108.572
$ws.Range("L83").Value = 37456
$ws.Range("M83").Value = -15116.572
$ws.Range("N83").Value = -47440
$ws.Range("H102").Value = 6116.4346
$ws.Range("I102").Value = 4642.5557
$ws.Range("K102").Value = 4642.5557
$ws.Range("M102").Value = -3020.5557
$ws.Range("H122").Value = 5245.7856
$ws.Range("I122").Value = 2264.6667
$ws.Range("J122").Value = 6058.8184
$ws.Range("K122").Value = 6794.000100000001
$ws.Range("L122").Value = 18176.4552
$ws.Range("M122").Value = -4344.000100000001
$ws.Range("N122").Value = -23076.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 37501250
$ws.Range("I4").Value = 33335000
$ws.Range("K4").Value = 33335000
$ws.Range("M4").Value = -33334887
$ws.Range("H25").Value = 10009199
$ws.Range("I25").Value = 16676667
$ws.Range("J25").Value = 7997.5
$ws.Range("K25").Value = 16676667
$ws.Range("L25").Value = 7997.5
$ws.Range("M25").Value = -16676437
$ws.Range("N25").Value = -8457.5
$ws.Range("H28").Value = 37501250
$ws.Range("I28").Value = 33335000
$ws.Range("K28").Value = 33335000
$ws.Range("M28").Value = -33334768
$ws.Range("H37").Value = 37501250
$ws.Range("I37").Value = 33335000
$ws.Range("K37").Value = 33335000
$ws.Range("M37").Value = -33334893
$ws.Range("H50").Value = 49499
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H81").Value = 29000
$ws.Range("J81").Value = 29000
$ws.Range("L81").Value = 29000
$ws.Range("N81").Value = -30996
$ws.Range("H84").Value = 29000
$ws.Range("J84").Value = 29000
$ws.Range("L84").Value = 87000
$ws.Range("N84").Value = -96984
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 8100
$ws.Range("J6").Value = 1800
$ws.Range("L6").Value = 1800
$ws.Range("N6").Value = -2030
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
